$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'polyester athletic pants men'
$ws.Range("A2").Value = 'compression for the knee'
$ws.Range("A3").Value = 'volleyball youth'
$ws.Range("A4").Value = 'knee sleeve with padding'
$ws.Range("A5").Value = 'kneeling pad exercise'
$ws.Range("A6").Value = 'baseball items'
$ws.Range("A7").Value = 'boys cold leggings'
$ws.Range("A8").Value = 'patella knee band'
$ws.Range("A9").Value = 'leg sleeves for basketball youth boys'
$ws.Range("A10").Value = 'mesh capri leggings'
$ws.Range("A11").Value = 'mens black compression pants'
$ws.Range("A12").Value = 'best knee pads construction'
$ws.Range("A13").Value = 'paintball pants small'
$ws.Range("A14").Value = 'below the knee'
$ws.Range("A15").Value = 'football pads'
$ws.Range("A16").Value = 'leg sleeves for basketball youth'
$ws.Range("A17").Value = 'knee pads for working'
$ws.Range("A18").Value = 'volleyball shorts longer length'
$ws.Range("A19").Value = 'kneepads construction'
$ws.Range("A20").Value = 'youth softball pants black'
$ws.Range("A21").Value = 'mountain skin pants'
$ws.Range("A22").Value = 'padded snowboarding shorts'
$ws.Range("A23").Value = 'volleyball kneepads'
$ws.Range("A24").Value = 'boys basketball gear'
$ws.Range("A25").Value = 'baseball leg guard'
$ws.Range("A26").Value = 'construction gel knee pads'
$ws.Range("A27").Value = 'girls youth softball pants black'
$ws.Range("A28").Value = 'running knee band'
$ws.Range("A29").Value = 'recovery pants'
$ws.Range("A30").Value = 'knee foam pad'
$ws.Range("A31").Value = 'lacrosse shorts girls'
$ws.Range("A32").Value = 'compression sleeve knee pads'
$ws.Range("A33").Value = 'compression running tights men'
$ws.Range("A34").Value = 'bjj shorts for men'
$ws.Range("A35").Value = 'football shorts for men'
$ws.Range("A36").Value = 'weight basketball'
$ws.Range("A37").Value = 'professional work knee pads'
$ws.Range("A38").Value = 'wrestling shorts'
$ws.Range("A39").Value = 'mesh basketball shorts for men'
$ws.Range("A40").Value = 'catchers gear leg guard'
$ws.Range("A41").Value = 'gel knee pads construction'
$ws.Range("A42").Value = 'acl knee'
$ws.Range("A43").Value = 'knee sleeve wrestling'
$ws.Range("A44").Value = 'spandex shorts men'
$ws.Range("A45").Value = 'sheer test'
$ws.Range("A46").Value = 'black leggings youth'
$ws.Range("A47").Value = 'calf sleeve padded'
$ws.Range("A48").Value = 'compression hip'
$ws.Range("A49").Value = 'knee sleeve for wrestling'
$ws.Range("A50").Value = 'compression knee for men'
$ws.Range("A51").Value = 'work knee pads'
$ws.Range("A52").Value = 'knees pads'
$ws.Range("A53").Value = 'medium youth baseball pants'
$ws.Range("A54").Value = 'knee compression sleeve with knee pad'
$ws.Range("A55").Value = 'knee pads with gel'
$ws.Range("A56").Value = 'football compression shorts youth'
$ws.Range("A57").Value = 'pants mountain'
$ws.Range("A58").Value = 'hex soccer'
$ws.Range("A59").Value = 'knee guards mountain biking'
$ws.Range("A60").Value = 'construction kneeling pad'
$ws.Range("A61").Value = 'leg guards softball'
$ws.Range("A62").Value = 'sort pants men'
$ws.Range("A63").Value = 'knee pad work'
$ws.Range("A64").Value = 'construction work knee pads'
$ws.Range("A65").Value = 'knee sleeves football'
$ws.Range("A66").Value = 'knee sleeve running men'
$ws.Range("A67").Value = 'squat pants men'
$ws.Range("A68").Value = 'down pants'
$ws.Range("A69").Value = 'adult baseball pants black'
$ws.Range("A70").Value = 'long basketball shorts'
$ws.Range("A71").Value = 'knee padding'
$ws.Range("A72").Value = 'shorts pad'
$ws.Range("A73").Value = 'leg sleeves for men football'
$ws.Range("A74").Value = 'black baseball pants youth'
$ws.Range("A75").Value = 'baseball shorts for boys'
$ws.Range("A76").Value = 'construction knee pads gel'
$ws.Range("A77").Value = 'mountain biking knee pads'
$ws.Range("A78").Value = 'mountain biking pads for men'
$ws.Range("A79").Value = 'work kneepads'
$ws.Range("A80").Value = 'construction knee pads for work'
$ws.Range("A81").Value = 'knee pads for works'
$ws.Range("A82").Value = 'compression knee sleeves for weightlifting'
$ws.Range("A83").Value = 'knee pads for work'
$ws.Range("A84").Value = 'compression shorts football'
$ws.Range("A85").Value = 'under shorts for men'
$ws.Range("A86").Value = 'soccer tights'
$ws.Range("A87").Value = 'knee pad floor'
$ws.Range("A88").Value = 'gel work knee pads'
$ws.Range("A89").Value = 'knee pads for kneeling'
$ws.Range("A90").Value = 'knees pads for work'
$ws.Range("A91").Value = 'biking capris'
$ws.Range("A92").Value = 'hip protector'
$ws.Range("A93").Value = 'volleyball kneepads women'
$ws.Range("A94").Value = 'compression shorts bjj'
$ws.Range("A95").Value = 'basketball sleeve youth boys'
$ws.Range("A96").Value = 'yoga tights'
$ws.Range("A97").Value = 'padded sliding shorts women'
$ws.Range("A98").Value = 'tight leggings'
$ws.Range("A99").Value = 'catchers gear women'
$ws.Range("A100").Value = 'softball gear for men'
